# 自动更新Excel文件 - Fri Nov 21 23:24:25 UTC 2025
# 重新计算"剩余"天数 (剩余 = 总天 - 自开始时间以来经过的天数)；
# 若已到期 (剩余 <= 0)，则以当天为新的开始时间重置周期 (剩余 = 总天)。

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 按行号记录新的 "剩余"(列E) 与 "开始时间"(列F) 取值，
# 对应把参考"今天"从 2025-11-21 推进到 2025-11-22 后的结果。
$updates = @{
    2 = @{ E = 9; F = 20251117 }
    3 = @{ E = 9; F = 20251117 }
    4 = @{ E = 9; F = 20251117 }
    5 = @{ E = 3; F = 20251115 }
    6 = @{ E = 9; F = 20251117 }
    7 = @{ E = 3; F = 20251115 }
    8 = @{ E = 9; F = 20251117 }
    9 = @{ E = 3; F = 20251115 }
    10 = @{ E = 2; F = 20251117 }
    11 = @{ E = 9; F = 20251117 }
    12 = @{ E = 3; F = 20251115 }
    13 = @{ E = 9; F = 20251117 }
    14 = @{ E = 9; F = 20251117 }
    15 = @{ E = 9; F = 20251117 }
    16 = @{ E = 7; F = 20251119 }
    17 = @{ E = 3; F = 20251115 }
    18 = @{ E = 6; F = 20251118 }
    19 = @{ E = 6; F = 20251118 }
    20 = @{ E = 6; F = 20251118 }
    21 = @{ E = 6; F = 20251118 }
    22 = @{ E = 3; F = 20251115 }
    23 = @{ E = 3; F = 20251115 }
    24 = @{ E = 3; F = 20251115 }
    25 = @{ E = 3; F = 20251115 }
    26 = @{ E = 3; F = 20251115 }
    27 = @{ E = 3; F = 20251118 }
    28 = @{ E = 6; F = 20251118 }
    29 = @{ E = 6; F = 20251118 }
    30 = @{ E = 6; F = 20251118 }
    31 = @{ E = 6; F = 20251118 }
    32 = @{ E = 6; F = 20251118 }
    33 = @{ E = 6; F = 20251118 }
    34 = @{ E = 6; F = 20251118 }
    35 = @{ E = 6; F = 20251118 }
    37 = @{ E = 6; F = 20251118 }
    38 = @{ E = 6; F = 20251118 }
    39 = @{ E = 6; F = 20251118 }
    40 = @{ E = 2; F = 20251117 }
    41 = @{ E = 2; F = 20251117 }
    42 = @{ E = 6; F = 20251118 }
    43 = @{ E = 3; F = 20251115 }
    44 = @{ E = 2; F = 20251117 }
    45 = @{ E = 3; F = 20251115 }
    46 = @{ E = 2; F = 20251117 }
    47 = @{ E = 6; F = 20251118 }
    48 = @{ E = 2; F = 20251117 }
    49 = @{ E = 3; F = 20251118 }
    50 = @{ E = 1; F = 20251113 }
    51 = @{ E = 1; F = 20251113 }
    52 = @{ E = 1; F = 20251113 }
    53 = @{ E = 1; F = 20251113 }
    54 = @{ E = 1; F = 20251113 }
    55 = @{ E = 1; F = 20251113 }
    56 = @{ E = 1; F = 20251113 }
    57 = @{ E = 1; F = 20251113 }
    58 = @{ E = 5; F = 20251117 }
    59 = @{ E = 5; F = 20251117 }
    60 = @{ E = 5; F = 20251117 }
    61 = @{ E = 3; F = 20251118 }
    62 = @{ E = 5; F = 20251117 }
    63 = @{ E = 5; F = 20251117 }
    64 = @{ E = 5; F = 20251117 }
    65 = @{ E = 6; F = 20251118 }
    66 = @{ E = 6; F = 20251118 }
    67 = @{ E = 6; F = 20251118 }
    68 = @{ E = 6; F = 20251118 }
    69 = @{ E = 6; F = 20251118 }
    70 = @{ E = 7; F = 20251119 }
    71 = @{ E = 7; F = 20251119 }
    72 = @{ E = 7; F = 20251119 }
    73 = @{ E = 7; F = 20251119 }
    74 = @{ E = 7; F = 20251119 }
    75 = @{ E = 7; F = 20251119 }
    76 = @{ E = 7; F = 20251119 }
    77 = @{ E = 10; F = 20251122 }
    78 = @{ E = 10; F = 20251122 }
    79 = @{ E = 10; F = 20251122 }
    80 = @{ E = 10; F = 20251122 }
    81 = @{ E = 10; F = 20251122 }
    82 = @{ E = 10; F = 20251122 }
    83 = @{ E = 10; F = 20251122 }
    84 = @{ E = 10; F = 20251122 }
    85 = @{ E = 10; F = 20251122 }
    86 = @{ E = 10; F = 20251122 }
    87 = @{ E = 2; F = 20251117 }
    88 = @{ E = 2; F = 20251117 }
    89 = @{ E = 2; F = 20251117 }
    90 = @{ E = 2; F = 20251117 }
    91 = @{ E = 3; F = 20251115 }
    92 = @{ E = 2; F = 20251117 }
    93 = @{ E = 10; F = 20251122 }
    94 = @{ E = 5; F = 20251120 }
    95 = @{ E = 9; F = 20251121 }
    96 = @{ E = 7; F = 20251119 }
    97 = @{ E = 7; F = 20251119 }
    98 = @{ E = 7; F = 20251119 }
    99 = @{ E = 7; F = 20251119 }

}

foreach ($rowNum in $updates.Keys) {
    $vals = $updates[$rowNum]
    $ws.Cells.Item($rowNum, 5).Value = $vals.E
    $ws.Cells.Item($rowNum, 6).Value = $vals.F
}
